$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 513.95
$ws.Range("I15").Value = 513.95
$ws.Range("K15").Value = 1541.85
$ws.Range("M15").Value = -1372.85
$ws.Range("H33").Value = 240.2
$ws.Range("I33").Value = 250.25
$ws.Range("K33").Value = 250.25
$ws.Range("M33").Value = -21.25
$ws.Range("H70").Value = 4224.385
$ws.Range("I70").Value = 2232.75
$ws.Range("J70").Value = 5109.5557
$ws.Range("K70").Value = 6698.25
$ws.Range("L70").Value = 15328.6671
$ws.Range("M70").Value = -6428.25
$ws.Range("N70").Value = -15868.6671
$ws.Range("H73").Value = 4224.385
$ws.Range("I73").Value = 2232.75
$ws.Range("J73").Value = 5109.5557
$ws.Range("K73").Value = 6698.25
$ws.Range("L73").Value = 15328.6671
$ws.Range("M73").Value = -5762.25
$ws.Range("N73").Value = -17200.6671
$ws.Range("H86").Value = 500
$ws.Range("I86").Value = 500
$ws.Range("K86").Value = 500
$ws.Range("M86").Value = 623
$ws.Range("H88").Value = 556.4286
$ws.Range("I88").Value = 580
$ws.Range("J88").Value = 497.5
$ws.Range("K88").Value = 580
$ws.Range("L88").Value = 497.5
$ws.Range("M88").Value = -174
$ws.Range("N88").Value = -1309.5
$ws.Range("H89").Value = 500
$ws.Range("I89").Value = 500
$ws.Range("K89").Value = 2500
$ws.Range("M89").Value = 3116
$ws.Range("H91").Value = 556.4286
$ws.Range("I91").Value = 580
$ws.Range("J91").Value = 497.5
$ws.Range("K91").Value = 580
$ws.Range("L91").Value = 497.5
$ws.Range("M91").Value = 824
$ws.Range("N91").Value = -3305.5
$ws.Range("H98").Value = 263.375
$ws.Range("I98").Value = 263.375
$ws.Range("K98").Value = 263.375
$ws.Range("M98").Value = 1234.625
$ws.Range("H122").Value = 263.375
$ws.Range("I122").Value = 263.375
$ws.Range("K122").Value = 790.125
$ws.Range("M122").Value = 1659.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 536.6786
$ws.Range("I32").Value = 536.6786
$ws.Range("K32").Value = 536.6786
$ws.Range("M32").Value = -249.6786
$ws.Range("H45").Value = 2123.1765
$ws.Range("I45").Value = 1259
$ws.Range("J45").Value = 3357.7144
$ws.Range("K45").Value = 1259
$ws.Range("L45").Value = 3357.7144
$ws.Range("M45").Value = -882
$ws.Range("N45").Value = -4111.7144
$ws.Range("H74").Value = 7082.3335
$ws.Range("I74").Value = 7082.3335
$ws.Range("K74").Value = 7082.3335
$ws.Range("M74").Value = -6208.3335
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 7082.3335
$ws.Range("I77").Value = 7082.3335
$ws.Range("K77").Value = 35411.6675
$ws.Range("M77").Value = -31043.6675
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H97").Value = 1790.6666
$ws.Range("I97").Value = 1798.8
$ws.Range("J97").Value = 1750
$ws.Range("K97").Value = 1798.8
$ws.Range("L97").Value = 1750
$ws.Range("M97").Value = -1302.8
$ws.Range("N97").Value = -2742

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3032.423
$ws.Range("I86").Value = 1502.15
$ws.Range("K86").Value = 1502.15
$ws.Range("M86").Value = -379.1500000000001
$ws.Range("H89").Value = 3032.423
$ws.Range("I89").Value = 1502.15
$ws.Range("K89").Value = 7510.75
$ws.Range("M89").Value = -1894.75
$ws.Range("H107").Value = 4825.25
$ws.Range("I107").Value = 3278.5
$ws.Range("K107").Value = 3278.5
$ws.Range("M107").Value = -1358.5
$ws.Range("H134").Value = 1592.7059
$ws.Range("I134").Value = 1379.75
$ws.Range("K134").Value = 4139.25
$ws.Range("M134").Value = -1604.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1992.4
$ws.Range("I132").Value = 1992.4
$ws.Range("K132").Value = 5977.200000000001
$ws.Range("M132").Value = -3447.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2945.4443
$ws.Range("I80").Value = 3203.5
$ws.Range("J80").Value = 2739
$ws.Range("K80").Value = 3203.5
$ws.Range("L80").Value = 2739
$ws.Range("M80").Value = -2205.5
$ws.Range("N80").Value = -4735
$ws.Range("H83").Value = 2945.4443
$ws.Range("I83").Value = 3203.5
$ws.Range("J83").Value = 2739
$ws.Range("K83").Value = 16017.5
$ws.Range("L83").Value = 13695
$ws.Range("M83").Value = -11025.5
$ws.Range("N83").Value = -23679
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 3835.95
$ws.Range("I113").Value = 2758.5715
$ws.Range("K113").Value = 2758.5715
$ws.Range("M113").Value = -588.5715
$ws.Range("H126").Value = 2256.5715
$ws.Range("I126").Value = 1949.5
$ws.Range("K126").Value = 5848.5
$ws.Range("M126").Value = -3378.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2750
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -888
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H82").Value = 4082.625
$ws.Range("I82").Value = 2372.2
$ws.Range("J82").Value = 6933.3335
$ws.Range("K82").Value = 2372.2
$ws.Range("L82").Value = 6933.3335
$ws.Range("M82").Value = -2011.2
$ws.Range("N82").Value = -7655.3335
$ws.Range("H85").Value = 4082.625
$ws.Range("I85").Value = 2372.2
$ws.Range("J85").Value = 6933.3335
$ws.Range("K85").Value = 2372.2
$ws.Range("L85").Value = 6933.3335
$ws.Range("M85").Value = -1124.2
$ws.Range("N85").Value = -9429.333500000001
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6235.7144
$ws.Range("I126").Value = 4485.7144
$ws.Range("J126").Value = 7985.7144
$ws.Range("K126").Value = 13457.1432
$ws.Range("L126").Value = 23957.1432
$ws.Range("M126").Value = -10987.1432
$ws.Range("N126").Value = -28897.1432
$ws.Range("H132").Value = 1669.125
$ws.Range("I132").Value = 951
$ws.Range("J132").Value = 2387.25
$ws.Range("K132").Value = 2853
$ws.Range("L132").Value = 7161.75
$ws.Range("M132").Value = -323
$ws.Range("N132").Value = -12221.75
$ws.Range("H135").Value = 69857.5
$ws.Range("J135").Value = 69857.5
$ws.Range("L135").Value = 69857.5
$ws.Range("N135").Value = -79997.5
$ws.Range("H136").Value = 2718.7666
$ws.Range("I136").Value = 1681.5714
$ws.Range("K136").Value = 5044.7142
$ws.Range("M136").Value = -2494.7142
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
